$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; the old row 25 (ImportMojio / Testing
# Required / Y) shifts down to become row 26.
$ws.Rows.Item(25).Insert()

# Give the freshly-inserted row the same formatting as the row above it
# (border/fill/font) without creating a duplicate style entry.
$ws.Range("A24:C24").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)

# New row 25: ImportMojio / Done / N
$ws.Range("A25").Value = "ImportMojio"
$ws.Range("B25").Value = "Done"
$ws.Range("C25").Value = "N"

# Row 24 description changes from "Testing Required" to "Done"
$ws.Range("B24").Value = "Done"

# Row 26 (the old row 25, now shifted down): ImportMojio/Testing Required/Y
# becomes ExportEvents / (blank) / Y
$ws.Range("A26").Value = "ExportEvents"
$ws.Range("B26").Value = ""

# Match the author's final selection/cursor position
$ws.Range("B19").Select()
